# Fix Training Data Issue (#48)
#
# The "Date" column (BF) on the single worksheet contained a malformed,
# constant placeholder value ("5-15-2007-08") for every data row instead
# of the actual game date. Correct it to the real ISO-style date
# "2008-05-15" for every data row (rows 2-31; row 1 is the "Date" header).
#
# The replacement text looks like a date, so a plain assignment would be
# auto-converted by Excel into a date serial number. Prefixing the value
# with a leading apostrophe forces Excel to keep it as literal text,
# which is what the source data expects (a text value, not a date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BF2:BF31").Value = "'2008-05-15"
